$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two obsolete 2025-07-09 rows (BEMOL S/A + MATHEUS SILVEIRA).
# Everything below shifts up by two rows.
$ws.Rows("2:3").Delete()

# Old row 4 (2025-07-11, MOCHILA PELUCIA STITCH) is now row 2; its date
# needs to move to 2025-07-15. Force text so Excel doesn't coerce it to a
# date serial, then drop the temporary formatting again.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2025-07-15"
$ws.Range("A2").ClearFormats()

# Old row 6 (BARALHO PLASTICO ...) is now row 4 - refreshed stats.
$ws.Range("G4").Value = -35
$ws.Range("H4").Value = 1.21
$ws.Range("I4").Value = 0.71

# Old row 8 (CARREGADOR USB-C ...) is now row 6 - refreshed stock delta.
$ws.Range("G6").Value = -94

# Insert a brand-new row for "MARMITA ELETRICA ONEX" ahead of the
# POWER BANK row (currently row 7, the former row 9), pushing it to row 8.
$ws.Rows("7:7").Insert()

$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "2025-07-23"
$ws.Range("A7").ClearFormats()

$ws.Range("B7").Value = 2
$ws.Range("C7").Value = "BEMOL S/A"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "386260"
$ws.Range("D7").ClearFormats()

$ws.Range("E7").Value = 47869
$ws.Range("F7").Value = "MARMITA ELÉTRICA ONEX"
$ws.Range("G7").Value = -16
$ws.Range("H7").Value = 1.09
$ws.Range("I7").Value = 0.3
